$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "fecha_inicial" / "fecha_final" data row with specific date values
# (previously both columns held "21"; now they hold distinct values "17" and "5").
# Leading apostrophe keeps them stored as text (shared strings) with the
# existing cell formatting/style untouched.
$ws.Range("C2").Value = "'17"
$ws.Range("D2").Value = "'5"

# Move the active selection to D9, matching the saved worksheet view state.
$null = $ws.Range("D9").Select()
